# Recalculated "DI generado" simulation output values for column A (rows 2-37).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2.1360139397911224
$ws.Range("A3").Value = 0.36872120936504243
$ws.Range("A4").Value = 0.50575891785357996
$ws.Range("A5").Value = 0.31365851114517557
$ws.Range("A6").Value = 0.26093017865784601
$ws.Range("A7").Value = 0.37867133240150952
$ws.Range("A8").Value = 0.068906038262757985
$ws.Range("A9").Value = 0.044182055587727809
$ws.Range("A10").Value = 0.20087840688674374
$ws.Range("A11").Value = 0.11720761929005896
$ws.Range("A12").Value = 0.19231661745942519
$ws.Range("A13").Value = 0.086437977771908592
$ws.Range("A14").Value = 0.11882932774353336
$ws.Range("A15").Value = 0.099169656590752053
$ws.Range("A16").Value = 0.023777124192886384
$ws.Range("A17").Value = 0.026667292763101903
$ws.Range("A18").Value = 0.05873889055617007
$ws.Range("A19").Value = 0.13647181440632178
$ws.Range("A20").Value = 0.041508201464209332
$ws.Range("A21").Value = 0.12547377449851252
$ws.Range("A22").Value = 0.051774328763595973
$ws.Range("A23").Value = 0.051176462020794011
$ws.Range("A24").Value = 0.023111647399819371
$ws.Range("A25").Value = 0.022175170233107087
$ws.Range("A26").Value = 0.03792953825870151
$ws.Range("A27").Value = 0.027297144788880987
$ws.Range("A28").Value = 0.025236278779534403
$ws.Range("A29").Value = 0.02649959749807251
$ws.Range("A30").Value = 0.030518950516186621
$ws.Range("A31").Value = 0.027731874171922563
$ws.Range("A32").Value = 0.025014689404894343
$ws.Range("A33").Value = 0.02660239922003271
$ws.Range("A34").Value = 0.026611369521752328
$ws.Range("A35").Value = 0.024557125901350126
$ws.Range("A36").Value = 0.035456452384477581
$ws.Range("A37").Value = 0.029053029017104062
